$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "25×39=" "76×46="
Replace-Text "44×64=" "53×67="
Replace-Text "59×19=" "39×94="
Replace-Text "87×52=" "50×64="
Replace-Text "31×99=" "51×51="
Replace-Text "28×11=" "76×67="
Replace-Text "83×89=" "50×67="
Replace-Text "12×19=" "89×69="
Replace-Text "28×33=" "70×70="
Replace-Text "44×32=" "64×12="
Replace-Text "79×28=" "26×90="
Replace-Text "83×94=" "26×45="
Replace-Text "37×33=" "54×68="
Replace-Text "19×38=" "68×31="
Replace-Text "43×55=" "11×95="
Replace-Text "33×46=" "36×79="
Replace-Text "55×65=" "63×72="
Replace-Text "65×63=" "28×17="
Replace-Text "63×15=" "81×40="
Replace-Text "61×56=" "12×93="
Replace-Text "18×58=" "37×14="
Replace-Text "83×21=" "76×69="
Replace-Text "67×65=" "70×98="
Replace-Text "25×29=" "34×40="
Replace-Text "46×99=" "99×30="
